$wb = $excel.ActiveWorkbook

# --- Rename existing "emojis" sheet to "emoj", insert a brand-new empty
# --- "emojis" sheet right after it (before "Draft"). ---
$emojSheet = $wb.Worksheets.Item("emojis")
$emojSheet.Name = "emoj"

$newEmojis = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $emojSheet)
$newEmojis.Name = "emojis"

# --- "emoj" sheet (formerly "emojis"): selection moves to the whole used
# --- range A1:A19, and the sheet gets unprotected + gets an explicit
# --- portrait page setup instead of sheetProtection. ---
$emojSheet.Unprotect()
$emojSheet.PageSetup.Orientation = 1
$emojSheet.Range("A1:A19").Select()

# --- "main" sheet: add a "motorcycle" entry and move the selection. ---
$mainSheet = $wb.Worksheets.Item("main")
$mainSheet.Range("F15").Value = "motorcycle"

# Make sure "main" ends up the active/selected sheet and tab again, and set
# its selection to B19 (was B22).
$mainSheet.Activate()
$mainSheet.Range("B19").Select()
